# Update "想去人数" (interest count) figures in the "展览" and "全部类型" sheets
# to reflect the latest scrape output.

$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibitions) - rows 5, 7, 10, 12
$wsExhibit = $wb.Worksheets.Item("展览")
$wsExhibit.Range("F5").Value = 15
$wsExhibit.Range("F7").Value = 2086
$wsExhibit.Range("F10").Value = 1181
$wsExhibit.Range("F12").Value = 1061

# Sheet "全部类型" (All types) - rows 5, 7, 11, 13
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F5").Value = 15
$wsAll.Range("F7").Value = 2086
$wsAll.Range("F11").Value = 1181
$wsAll.Range("F13").Value = 1061
